$wb = $excel.ActiveWorkbook

# --- Schedule sheet ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 789.1812434999999
$schedule.Range("F2").Value = 13.04863167162698

# --- Detailed sheet ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B18").Value = 0.7

$detailed.Range("B19").Value = 36.06011
$detailed.Range("C19").Value = "historical"

$detailed.Range("B20").Value = 36.06011
$detailed.Range("C20").Value = "historical"

$detailed.Range("B22").Value = 5.98325

$detailed.Range("B24").Value = 32.34641

$detailed.Range("B25").Value = 36.06046

$detailed.Range("B26").Value = 35.88

$detailed.Range("B27").Value = 32.51579

$detailed.Range("B28").Value = 36.06045

$detailed.Range("B29").Value = 36.06029

$detailed.Range("B30").Value = 23.63823

$detailed.Range("B31").Value = 23.5144

$detailed.Range("B32").Value = 24.28971

$detailed.Range("B33").Value = 36.061

$detailed.Range("B34").Value = 33.26822

$detailed.Range("B35").Value = 23.74211

$detailed.Range("B36").Value = -0.12629

$detailed.Range("B37").Value = -3.02299

$detailed.Range("B38").Value = -2.9204

$detailed.Range("B39").Value = -2.90758

$detailed.Range("B40").Value = 4.5258

$detailed.Range("B41").Value = 32.90602

$detailed.Range("B42").Value = 33.03132

$detailed.Range("B43").Value = 9.53241

$detailed.Range("B44").Value = 19.51387

$detailed.Range("B45").Value = 36.0601

$detailed.Range("B47").Value = 57.3
